$d = $word.ActiveDocument

# The outline has two "Invariant tracking" bullets (one under Ch. 1, one
# under Ch. 6 "Simulations"). We want the Ch. 6 one, which is the bullet
# immediately followed by "Frequency map analysis [old research update]".
$search = $d.Content
$target = $null
while ($search.Find.Execute("Invariant tracking", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $matchEnd = $search.End
    $after = $d.Range($matchEnd, $matchEnd + 60)
    if ($after.Text.Contains("Frequency map analysis")) {
        $target = $d.Range($search.Start, $matchEnd)
        break
    }
    $search = $d.Range($matchEnd, $d.Content.End)
}

$insertAt = $target.End

# Remove the existing "_GoBack" bookmark -- it currently sits collapsed at
# the end of the "Tune scan calibration (Robust Method)" paragraph (the
# spot of the previous paste/edit).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create "_GoBack" collapsed at the end of the "Invariant tracking"
# paragraph (the new paste location), right before that paragraph's mark.
# A bookmark can't be added directly on a zero-length range that sits on
# a paragraph-mark boundary, so: insert a throwaway character, bookmark
# across it, then delete the character -- the bookmark collapses in place
# exactly where we want it, leaving the surrounding text untouched.
$marker = $d.Range($insertAt, $insertAt)
$marker.InsertAfter("X")
$markerSpan = $d.Range($insertAt, $insertAt + 1)
$d.Bookmarks.Add("_GoBack", $markerSpan)
$d.Range($insertAt, $insertAt + 1).Delete()
